$d = $word.ActiveDocument

# Locate the paragraph that must be preserved ("LOB1205: Ecologia Básica
# (Requisito)") and the last paragraph that must be removed ("© 2020 ...").
# Everything from right after the first paragraph's mark through the end
# of the second paragraph's mark (i.e. the blank paragraph, the
# "Ver no Jupiter..." paragraph and the "© 2020..." paragraph) is deleted,
# leaving the trailing blank paragraph / page-break paragraph untouched.

$keepEnd = $d.Content.Duplicate
[void]$keepEnd.Find.Execute("LOB1205: Ecologia Básica (Requisito)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$removeEnd = $d.Content.Duplicate
[void]$removeEnd.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$junk = $d.Range($keepEnd.End + 1, $removeEnd.End + 1)
$junk.Delete()
